$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.309.24'
$ws.Range("E2").Value = '  -1.05%  '
$ws.Range("D3").Value = '1.561.45'
$ws.Range("E3").Value = '  -1.00%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.000'
$ws.Range("E5").Value = '  -0.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '289.82'
$ws.Range("E6").Value = '  +0.23%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3772'
$ws.Range("E7").Value = '  +1.90%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3296'
$ws.Range("E8").Value = '  -1.58%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.70'
$ws.Range("E10").Value = '  +0.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07400'
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.37'
$ws.Range("E13").Value = '  -3.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.887'
$ws.Range("E14").Value = '  -1.93%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.781'
$ws.Range("E15").Value = '  -2.45%  '
$ws.Range("D16").Value = '1.568.62'
$ws.Range("E16").Value = '  -0.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001082'
$ws.Range("E17").Value = '  -3.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06660'
$ws.Range("E18").Value = '  -1.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '86.53'
$ws.Range("E19").Value = '  -2.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.441'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9999'
$ws.Range("E21").Value = '  -0.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.24'
$ws.Range("E22").Value = '  -2.10%  '
$ws.Range("E23").Value = '  -3.33%  '
$ws.Range("D24").Value = '22.305.70'
$ws.Range("E24").Value = '  -1.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.298'
$ws.Range("E25").Value = '  -4.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.600'
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.25'
$ws.Range("E27").Value = '  -1.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.35'
$ws.Range("E28").Value = '  -1.84%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.935'
$ws.Range("E29").Value = '  -1.61%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.38'
$ws.Range("E30").Value = '  -0.81%  '
$ws.Range("D31").Value = '1.742.80'
$ws.Range("E31").Value = '  -0.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.088'
$ws.Range("E32").Value = '  +2.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.960'
$ws.Range("E33").Value = '  -3.63%  '
$ws.Range("E34").Value = '  -4.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.490'
$ws.Range("E35").Value = '  -1.88%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08208'
$ws.Range("E36").Value = '  -1.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02369'
$ws.Range("E37").Value = '  -3.61%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06371'
$ws.Range("E38").Value = '  -0.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.391'
$ws.Range("E39").Value = '  -0.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2173'
$ws.Range("E40").Value = '  -5.71%  '
$ws.Range("E41").Value = '  -3.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.09'
$ws.Range("E42").Value = '  -2.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6107'
$ws.Range("E43").Value = '  -3.88%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9997'
$ws.Range("E44").Value = '  -0.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.91'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.767'
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5933'
$ws.Range("E47").Value = '  -4.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.993'
$ws.Range("E48").Value = '  -3.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '123.71'
$ws.Range("E49").Value = '  -0.90%  '
$ws.Range("E50").Value = '  -2.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07107'
